# Sync attendance_reports: reorder "Recorded By" (column G) names so that
# "System" / "system" and any backup/admin identifiers are listed before the
# primary account, matching the canonical ordering from the main repo.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2   = "system, System, backup@backdoor.com"
    4   = "System, backup@backdoor.com"
    5   = "System, backup@backdoor.com"
    8   = "System, backup@backdoor.com"
    11  = "System, dnasr281@gmail.com"
    17  = "System, dnasr281@gmail.com"
    28  = "system, System, backup@backdoor.com"
    30  = "System, backup@backdoor.com"
    31  = "System, backup@backdoor.com"
    34  = "System, backup@backdoor.com"
    37  = "System, dnasr281@gmail.com"
    43  = "System, dnasr281@gmail.com"
    54  = "system, System, backup@backdoor.com"
    56  = "System, backup@backdoor.com"
    57  = "System, backup@backdoor.com"
    60  = "System, backup@backdoor.com"
    63  = "System, dnasr281@gmail.com"
    69  = "System, dnasr281@gmail.com"
    80  = "System, backup@backdoor.com"
    81  = "System, backup@backdoor.com"
    82  = "System, backup@backdoor.com"
    93  = "System, dnasr281@gmail.com"
    94  = "System, dnasr281@gmail.com"
    96  = "System, dnasr281@gmail.com"
    106 = "System, backup@backdoor.com"
    107 = "System, backup@backdoor.com"
    108 = "System, backup@backdoor.com"
    119 = "System, dnasr281@gmail.com"
    120 = "System, dnasr281@gmail.com"
    122 = "System, dnasr281@gmail.com"
    132 = "System, backup@backdoor.com"
    133 = "System, backup@backdoor.com"
    134 = "System, backup@backdoor.com"
    145 = "System, dnasr281@gmail.com"
    146 = "System, dnasr281@gmail.com"
    148 = "System, dnasr281@gmail.com"
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
